$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.900.91"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.635.77"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.04%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "211.74"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -0.53%  "
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.06%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.34"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +0.67%  "
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -0.19%  "
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -0.08%  "
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.868.29"
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.645.45"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +0.67%  "
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -0.67%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.562"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -0.79%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.39"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +0.39%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "27.918.49"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "228.85"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -0.53%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.69"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +2.57%  "
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -0.09%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.09"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -2.16%  "
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +0.34%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "155.35"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +1.71%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.91"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -0.60%  "
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +0.19%  "
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -0.37%  "
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +0.32%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0482"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +0.06%  "
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +1.32%  "
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +1.35%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.398.22"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +3.64%  "
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +0.18%  "
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -0.78%  "
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +0.76%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.560"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -2.08%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -1.65%  "
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +2.82%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "66.09"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -1.09%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "5.47"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -0.74%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.776.99"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -2.76%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "88.77"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +1.26%  "
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +2.66%  "
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -0.32%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +1.41%  "
$c.Style = "Normal"
